$d = $word.ActiveDocument

# Locate the target paragraph: the one ending in
# "... محصولات فوق تخصصی ... می رسد." by searching for a distinctive
# substring of its (only) run, then expanding to the enclosing paragraph.
$findRange = $d.Content
$null = $findRange.Find.Execute("به علاوه محصولات فوق تخصصی", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target = $findRange.Paragraphs(1)

# Index of the target paragraph, so we can re-fetch it (by position) after
# the document's paragraph collection shifts below.
$targetIndex = $target.Index

# 1) Split off an empty paragraph immediately BEFORE the target paragraph.
#    Word keeps the target paragraph's own mark/properties (paraId, rsids,
#    pPr, ...) on the (second) paragraph that keeps its original position;
#    the brand-new (first / inserted) paragraph is the "empty" one.
$startOfTarget = $target.Range
$startOfTarget.Collapse(1)
$null = $startOfTarget.InsertParagraphBefore()

# After the insert, the document gained one paragraph: the new empty one
# sits where the target used to be (index $targetIndex), and the original
# target (with all of its own content + attributes) is now one slot later.
$newPara = $d.Paragraphs($targetIndex)
$target  = $d.Paragraphs($targetIndex + 1)

# 2) Move the target's original run text into the new (currently empty)
#    paragraph, so it ends up first in reading order once merged back.
$targetRange = $target.Range
$targetRange.MoveEnd(1, -1) | Out-Null
$originalText = $targetRange.Text

$newRange = $newPara.Range
$newRange.MoveEnd(1, -1) | Out-Null
$newRange.Text = $originalText

# 3) Replace the (still-attributed) target paragraph's own text with just
#    the new trailing run: a single space, matching the added run's text.
$targetRange2 = $target.Range
$targetRange2.MoveEnd(1, -1) | Out-Null
$targetRange2.Text = " "

# 4) Delete the new (first) paragraph's own mark. This merges its run
#    (the original sentence) forward into the target paragraph, which
#    keeps its own paragraph mark / properties / paraId as the survivor -
#    yielding one paragraph with two runs: the untouched original run,
#    then the new single-space run, exactly as in the target diff.
$mark = $d.Range($newPara.Range.End - 1, $newPara.Range.End)
$mark.Delete()
